$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.255269829918575,0.1762097630465576,0.1273133655188019,0.1219173530157218,1.582352370776711,0.9683392973091784,0.1535881532751162,0.3136387744262805,0.3075114984643506,4.051546183316816),
    @(1.151776496250534,0.1600266875675231,0.1269731100013587,0.1227890251942423,1.59412817672257,0.9813925671496264,0.154557706433458,0.3096592936175284,0.2905356750174946,4.086702161348541),
    @(1.088267293465037,0.1500344830695042,0.1267937272366275,0.1233555329943685,1.602255175343437,0.9899683364321525,0.1551909129103066,0.3073118794476102,0.2801605470758801,4.11076440573261),
    @(1.062397676822172,0.1459488057901979,0.1267280878340564,0.1235942730178112,1.605792459628852,0.9936040380864739,0.155458499276774,0.3063795564470766,0.2759450119318956,4.121192419656296),
    @(1.058102758662869,0.1452695574817255,0.1267176400490904,0.1236343921565159,1.606393440824093,0.9942162581111802,0.155503509153851,0.306226214344143,0.275245784734885,4.122961571078434),
    @(1.087918360221352,0.1499794376407806,0.1267928117479897,0.1233587207878399,1.602301967483456,0.990016797934409,0.155194482980832,0.3072992073872172,0.280103644206811,4.110902521607443),
    @(1.219578884434839,0.1706415448983876,0.1271899351611481,0.1222114179346296,1.586226581612117,0.9727235898947555,0.1539146017823558,0.3122467980090917,0.3016484377176383,4.063153856397236),
    @(1.477984808977794,0.2107090024374827,0.1282016814550246,0.1202092679808815,1.561817346997103,0.9432656104050565,0.1517045478674941,0.3227059155216949,0.3442675539321058,3.989179689005539),
    @(1.667899729086969,0.2398628693835008,0.1290853826186265,0.1188884091111024,1.548222672725061,0.9243412458260103,0.1502623289237857,0.3308462316968246,0.3757927740186133,3.94683564734774),
    @(1.754297375352792,0.2530624431406068,0.129517541250955,0.1183199238024888,1.542980650693423,0.9163230144484267,0.1496453755647469,0.334647412772739,0.3901782338861466,3.9301832282639),
    @(1.787013009736484,0.258051566137425,0.1296854940006753,0.1181092955184306,1.541131153565281,0.913371714074195,0.1494173569160058,0.336100817927445,0.3956317567994745,3.924253051532531),
    @(1.779967199694397,0.2569774855626861,0.1296491314901331,0.118154451692384,1.541523446828108,0.9140035459651727,0.1494662156554654,0.335787181370037,0.394456978783019,3.925513503966755),
    @(1.756988946819718,0.2534730885276133,0.1295312727619304,0.1183025022641822,1.542825774707921,0.9160785044452879,0.1496265040337388,0.3347667058316262,0.3906267789263822,3.929687816502508),
    @(1.742913880367894,0.2513253294294486,0.1294596403711523,0.1183937920613136,1.543641141182796,0.9173605530377742,0.1497254151429033,0.3341434519743842,0.3882814504083925,3.932293646929594),
    @(1.662253364477863,0.2389989631050753,0.1290577440684189,0.1189262115519583,1.548584215939414,0.9248771488939731,0.1503034339890803,0.3305997788699386,0.3748535167682832,3.947976478545485),
    @(1.612770512264206,0.2314209021700719,0.1288188944180675,0.1192611187980455,1.551858001861106,0.9296396700124596,0.1506680369693711,0.328450887523374,0.3666270678350259,3.958266184863533),
    @(1.584309811643664,0.2270563187595371,0.1286843550512913,0.1194567970281923,1.553829684505786,0.9324345257346742,0.1508814300586501,0.3272241442760873,0.36189963843794,3.964430202551938),
    @(1.574673656451068,0.2255775432543032,0.1286392912263565,0.1195235742370264,1.554512493118146,0.9333903583951724,0.1509543143929406,0.326810381354079,0.3602997447643403,3.96655941143689),
    @(1.618038005682308,0.2322282104565545,0.1288440266168251,0.1192251519310483,1.551500322588609,0.9291269384477765,0.1506288432971825,0.3286786850688941,0.3675023544466924,3.957145401185301),
    @(1.763738266571579,0.2545026684907725,0.1295657741908727,0.1182588902556403,1.542439570480617,0.9154667302931152,0.1495792713821835,0.3350660658260125,0.391751639618974,3.928451519232596),
    @(1.858953476297529,0.2690061687759169,0.1300625477051511,0.1176544503777284,1.537307913591249,0.9070346637576918,0.1489259965597522,0.3393220178546557,0.4076351264063334,3.911888653490905),
    @(1.808136717430443,0.2612704166854485,0.1297951273407634,0.1179745783550263,1.539974469009096,0.9114896267982608,0.14927167689663,0.3370431274246499,0.3991547018611001,3.920528026319658),
    @(1.61565661023235,0.2318632508307985,0.1288326556931096,0.1192414027837674,1.551661750480015,0.9293585674073128,0.1506465509866182,0.3285756707427652,0.3671066309102358,3.957651333982852),
    @(1.408062899650133,0.1999188834311667,0.1279032070290782,0.1207244759200994,1.567658865321242,0.9507576211714124,0.1522704618071167,0.3197959785234161,0.3326996615312225,4.007084861261262)
)

$cols = @("B","C","D","E","F","I","J","L","M","O")

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $cellRef = "$($cols[$j])$rowNum"
        $ws.Range($cellRef).Value = $rowVals[$j]
    }
}
